$wb = $excel.ActiveWorkbook

# Rename "SEM NOTA FISCAL" -> "SEM_NOTA_FISCAL"
$wsSem = $wb.Worksheets.Item("SEM NOTA FISCAL")
$wsSem.Name = "SEM_NOTA_FISCAL"

# EXECUTADO 2019 - scroll up one row (topLeftCell A38 -> A37); selection stays I69
$ws2019 = $wb.Worksheets.Item("EXECUTADO 2019")
$ws2019.Activate()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws2019.Range("I69").Select()

# EXECUTADO 2020 - selection moves from A33 to A32, scroll reset (topLeftCell cleared)
$ws2020 = $wb.Worksheets.Item("EXECUTADO 2020")
$ws2020.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws2020.Range("A32").Select()

# SEM_NOTA_FISCAL - selection moves from G29 to D28 and becomes the active tab
$wsSem.Activate()
$wsSem.Range("D28").Select()
